$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'290.68"
$ws.Range("E2").Value = "'-4.16%"
$ws.Range("D3").Value = "'30.84"
$ws.Range("E3").Value = "'-6.81%"
$ws.Range("D4").Value = "'4.945"
$ws.Range("E4").Value = "'0.06%"
$ws.Range("D5").Value = "'0.07206"
$ws.Range("E5").Value = "'-8.04%"
$ws.Range("D6").Value = "'1.790"
$ws.Range("E6").Value = "'-10.93%"
$ws.Range("D7").Value = "'7.664"
$ws.Range("E7").Value = "'-2.20%"
$ws.Range("D8").Value = "'3.750"
$ws.Range("E8").Value = "'-1.51%"
$ws.Range("D9").Value = "'0.8958"
$ws.Range("E9").Value = "'-3.04%"
$ws.Range("D10").Value = "'0.1660"
$ws.Range("E10").Value = "'-5.56%"
$ws.Range("D11").Value = "'0.07734"
$ws.Range("E11").Value = "'-1.48%"
$ws.Range("D12").Value = "'0.07987"
$ws.Range("E12").Value = "'-8.03%"
$ws.Range("D13").Value = "'0.03067"
$ws.Range("E13").Value = "'-2.35%"
$ws.Range("E14").Value = "'-0.24%"
$ws.Range("D15").Value = "'0.001502"
$ws.Range("E15").Value = "'-0.71%"
$ws.Range("D16").Value = "'0.005667"
$ws.Range("E16").Value = "'-3.76%"
$ws.Range("E17").Value = "'0.26%"
$ws.Range("D18").Value = "'2.081"
$ws.Range("E18").Value = "'-3.43%"
$ws.Range("D19").Value = "'0.3279"
$ws.Range("E19").Value = "'-0.88%"
$ws.Range("D20").Value = "'0.1298"
$ws.Range("E20").Value = "'-1.55%"
$ws.Range("D21").Value = "'4.051"
$ws.Range("E21").Value = "'-6.11%"
$ws.Range("D22").Value = "'0.2101"
$ws.Range("E22").Value = "'5.61%"
$ws.Range("D23").Value = "'0.04510"
$ws.Range("E23").Value = "'-1.18%"
$ws.Range("D24").Value = "'0.001214"
$ws.Range("E24").Value = "'-0.87%"
$ws.Range("D25").Value = "'0.004006"
$ws.Range("E25").Value = "'-9.89%"
$ws.Range("D26").Value = "'0.0001251"
$ws.Range("E26").Value = "'0.09%"
$ws.Range("E39").Value = "'-8.02%"
$ws.Range("D40").Value = "'0.04377"
$ws.Range("E40").Value = "'-8.58%"
$ws.Range("D41").Value = "'0.007399"
$ws.Range("E41").Value = "'-1.28%"
$ws.Range("E42").Value = "'-3.86%"
$ws.Range("D43").Value = "'0.007676"
$ws.Range("D44").Value = "'0.002062"
$ws.Range("E44").Value = "'-11.89%"
$ws.Range("D45").Value = "'0.009211"
$ws.Range("E45").Value = "'-12.83%"
$ws.Range("D46").Value = "'0.00005915"
$ws.Range("E46").Value = "'-5.27%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.09%"
$ws.Range("E48").Value = "'172.73%"
$ws.Range("D49").Value = "'0.003002"
$ws.Range("E49").Value = "'-3.20%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.09%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'0.09%"
